$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 516.8333
$ws.Range("I6").Value = 520.2
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 1560.6
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -1448.6
$ws.Range("N6").Value = -1724

$ws.Range("H15").Value = 331.92
$ws.Range("I15").Value = 331.92
$ws.Range("K15").Value = 995.76
$ws.Range("M15").Value = -826.76

$ws.Range("H125").Value = 835.6111
$ws.Range("I125").Value = 465.2
$ws.Range("J125").Value = 978.0769
$ws.Range("K125").Value = 4186.8
$ws.Range("L125").Value = 8802.6921
$ws.Range("M125").Value = -1726.8
$ws.Range("N125").Value = -13722.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1885.8096
$ws.Range("I2").Value = 792.46155
$ws.Range("K2").Value = 792.46155
$ws.Range("M2").Value = -679.46155

$ws.Range("H61").Value = 3613.1428
$ws.Range("I61").Value = 3613.1428
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3613.1428
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3401.1428
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 2012.6666
$ws.Range("J74").Value = 2155.5557
$ws.Range("L74").Value = 2155.5557
$ws.Range("N74").Value = -3903.5557

$ws.Range("H77").Value = 2012.6666
$ws.Range("J77").Value = 2155.5557
$ws.Range("L77").Value = 10777.7785
$ws.Range("N77").Value = -19513.7785

$ws.Range("H116").Value = 1885.8096
$ws.Range("I116").Value = 792.46155
$ws.Range("K116").Value = 792.46155
$ws.Range("M116").Value = 1501.53845

$ws.Range("H136").Value = 3613.1428
$ws.Range("I136").Value = 3613.1428
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10839.4284
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8289.428400000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1885.8096
$ws.Range("I3").Value = 792.46155
$ws.Range("K3").Value = 792.46155
$ws.Range("M3").Value = -678.46155

$ws.Range("H105").Value = 1631.081
$ws.Range("I105").Value = 1547.0588
$ws.Range("J105").Value = 1702.5
$ws.Range("K105").Value = 1547.0588
$ws.Range("L105").Value = 1702.5
$ws.Range("M105").Value = 199.9412
$ws.Range("N105").Value = -5196.5

$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -330
$ws.Range("N12").ClearContents()

$ws.Range("H19").Value = 681.8182
$ws.Range("I19").Value = 233.33333
$ws.Range("J19").Value = 2700
$ws.Range("K19").Value = 233.33333
$ws.Range("L19").Value = 2700
$ws.Range("M19").Value = -63.33332999999999
$ws.Range("N19").Value = -3040

$ws.Range("H24").Value = 681.8182
$ws.Range("I24").Value = 233.33333
$ws.Range("J24").Value = 2700
$ws.Range("K24").Value = 233.33333
$ws.Range("L24").Value = 2700
$ws.Range("M24").Value = -63.33332999999999
$ws.Range("N24").Value = -3040

$ws.Range("H58").Value = 1202.2916
$ws.Range("I58").Value = 962.75
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 962.75
$ws.Range("L58").Value = 2400
$ws.Range("M58").Value = -759.75
$ws.Range("N58").Value = -2806

$ws.Range("H136").Value = 1202.2916
$ws.Range("I136").Value = 962.75
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 2888.25
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -338.25
$ws.Range("N136").Value = -12300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1428.5
$ws.Range("I63").Value = 1250
$ws.Range("J63").Value = 1607
$ws.Range("K63").Value = 3750
$ws.Range("L63").Value = 4821
$ws.Range("M63").Value = -3001
$ws.Range("N63").Value = -6319

$ws.Range("H66").Value = 1428.5
$ws.Range("I66").Value = 1250
$ws.Range("J66").Value = 1607
$ws.Range("K66").Value = 11250
$ws.Range("L66").Value = 14463
$ws.Range("M66").Value = -7506
$ws.Range("N66").Value = -21951

$ws.Range("H94").Value = 6500
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H131").Value = 6098429
$ws.Range("I131").Value = 1175
$ws.Range("J131").Value = 6757592
$ws.Range("K131").Value = 3525
$ws.Range("L131").Value = 20272776
$ws.Range("M131").Value = 1515
$ws.Range("N131").Value = -20282856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2543.739
$ws.Range("I80").Value = 1885.7142
$ws.Range("J80").Value = 3567.3333
$ws.Range("K80").Value = 1885.7142
$ws.Range("L80").Value = 3567.3333
$ws.Range("M80").Value = -887.7141999999999
$ws.Range("N80").Value = -5563.3333

$ws.Range("H83").Value = 2543.739
$ws.Range("I83").Value = 1885.7142
$ws.Range("J83").Value = 3567.3333
$ws.Range("K83").Value = 9428.571
$ws.Range("L83").Value = 17836.6665
$ws.Range("M83").Value = -4436.571
$ws.Range("N83").Value = -27820.6665

$ws.Range("H102").Value = 2453.158
$ws.Range("I102").Value = 2350.6667
$ws.Range("J102").Value = 2837.5
$ws.Range("K102").Value = 2350.6667
$ws.Range("L102").Value = 2837.5
$ws.Range("M102").Value = -728.6667000000002
$ws.Range("N102").Value = -6081.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 10003550
$ws.Range("I14").Value = 20002100
$ws.Range("J14").Value = 4999.5
$ws.Range("K14").Value = 20002100
$ws.Range("L14").Value = 4999.5
$ws.Range("M14").Value = -20001928
$ws.Range("N14").Value = -5343.5

$ws.Range("H40").Value = 2726.7646
$ws.Range("I40").Value = 2667.111
$ws.Range("J40").Value = 2793.875
$ws.Range("K40").Value = 2667.111
$ws.Range("L40").Value = 2793.875
$ws.Range("M40").Value = -2531.111
$ws.Range("N40").Value = -3065.875

$ws.Range("H68").Value = 1965.65
$ws.Range("I68").Value = 1890.6552
$ws.Range("J68").Value = 2163.3635
$ws.Range("K68").Value = 1890.6552
$ws.Range("L68").Value = 2163.3635
$ws.Range("M68").Value = -1141.6552
$ws.Range("N68").Value = -3661.3635

$ws.Range("H71").Value = 1965.65
$ws.Range("I71").Value = 1890.6552
$ws.Range("J71").Value = 2163.3635
$ws.Range("K71").Value = 9453.276
$ws.Range("L71").Value = 10816.8175
$ws.Range("M71").Value = -5709.276
$ws.Range("N71").Value = -18304.8175

$ws.Range("H122").Value = 2616.3333
$ws.Range("I122").Value = 2724.5
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 8173.5
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -5723.5
$ws.Range("N122").Value = -12100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H96").Value = 86070.836
$ws.Range("J96").Value = 114300
$ws.Range("L96").Value = 114300
$ws.Range("N96").Value = -117046

$ws.Range("H122").Value = 1809
$ws.Range("I122").Value = 1718.7693
$ws.Range("K122").Value = 5156.3079
$ws.Range("M122").Value = -2706.3079

$ws.Range("H136").Value = 876.06665
$ws.Range("I136").Value = 778.4167
$ws.Range("J136").Value = 1266.6666
$ws.Range("K136").Value = 2335.2501
$ws.Range("L136").Value = 3799.9998
$ws.Range("M136").Value = 214.7498999999998
$ws.Range("N136").Value = -8899.9998
